$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208, shifting rows 208:221 down to 209:222
$ws.Rows.Item(208).Insert()

$ws.Cells.Item(208, 1).Value = 3
$ws.Cells.Item(208, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(208, 3).Value = "Coquimbo"
$ws.Cells.Item(208, 4).Value = 44585
$ws.Cells.Item(208, 5).Value = 5
$ws.Cells.Item(208, 6).Value = 100112001
$ws.Cells.Item(208, 7).Value = "Berenjena"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 80
$ws.Cells.Item(208, 11).Value = 9500
$ws.Cells.Item(208, 12).Value = 10000
$ws.Cells.Item(208, 13).Value = 9750
$ws.Cells.Item(208, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(208, 15).Value = "Región Metropolitana"
$ws.Cells.Item(208, 16).Value = 162
$ws.Cells.Item(208, 17).Value = 60
$ws.Cells.Item(208, 18).Value = "Hortaliza"
